$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195984482765198
$ws.Range("B1").Value = 2.105360269546509
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.142131805419922
$ws.Range("E1").Value = 1.212249875068665
